$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4-6 down to 5-7
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = "c1243957"
$ws.Range("B4").Value = @"
Cursist ID:  c1243957 

 Gebruik de data in https://feb.kuleuven.be/public/U0004359/datayXzkle.txt 
 De vragen voor deze taak staan hieronder vermeld. 


 V1: Voer een regressie analyse uit met als afhankelijke variabele Y1 en als onafhankelijke variabelen X1 en X2 (zonder interactieterm). Geef de parameterschatting die hoort bij variabele X2. 

 V2: Voer een regressie analyse uit met als afhankelijke variabele Y1, en onafhankelijke variabelen X2 en X4. Neem ook de interactieterm op. Geef de p-waarde die bij de parameter van de interactieterm hoort. 

 V3: Voer een regressie analyse uit met als afhankelijke variabele Y3, en onafhankelijke variabelen X1, X2 en X3 (zonder interactietermen). Geef de proportie verklaarde variantie. 


 Vergeet kommagetallen niet af te ronden op 3 decimalen.
"@

$ws.Range("A5").Value = "q0762379"
$ws.Range("B5").Value = @"
Student ID:  q0762379 

 Use the data in https://feb.kuleuven.be/public/U0004359/datahPQFvR.txt 
 The questions for this task are listed below. 


 Q1: Perform a regression analysis with dependent variable Y2 and independent variables X1 and X3 (without the interaction term). Give the estimate for the parameter corresponding to variable X3. 

 Q2: Perform a regression analysis with dependent variable Y3 and independent variables X3 and X4. Take the interaction into account. Give the p-value of the parameter corresponding to the interaction term. 

 Q3: Perform a regression analysis with dependent variable Y3 and independent variables X1, X2 and X3 (without the interaction terms). Give the proportion of explained variance. 


 Don't forget to round decimals to three digits.
"@

$ws.Range("A6").Value = "q1371623"
$ws.Range("B6").Value = @"
Student ID:  q1371623 

 Use the data in https://feb.kuleuven.be/public/U0004359/datalDLHdt.txt 
 The questions for this task are listed below. 


 Q1: Perform a regression analysis with dependent variable Y1 and independent variables X1 and X2 (without the interaction term). Give the estimate for the parameter corresponding to variable X2. 

 Q2: Perform a regression analysis with dependent variable Y3 and independent variables X3 and X4. Take the interaction into account. Give the p-value of the parameter corresponding to the interaction term. 

 Q3: Perform a regression analysis with dependent variable Y3 and independent variables X1, X2 and X3 (without the interaction terms). Give the proportion of explained variance. 


 Don't forget to round decimals to three digits.
"@

$ws.Range("A7").Value = "q1411379"
$ws.Range("B7").Value = @"
Student ID:  q1411379 

 Use the data in https://feb.kuleuven.be/public/U0004359/dataPrlRAD.txt 
 The questions for this task are listed below. 


 Q1: Perform a regression analysis with dependent variable Y1 and independent variables X1 and X2 (without the interaction term). Give the estimate for the parameter corresponding to variable X2. 

 Q2: Perform a regression analysis with dependent variable Y1 and independent variables X2 and X4. Take the interaction into account. Give the p-value of the parameter corresponding to the interaction term. 

 Q3: Perform a regression analysis with dependent variable Y1 and independent variables X2, X3 and X4 (without the interaction terms). Give the proportion of explained variance. 


 Don't forget to round decimals to three digits.
"@

$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(7).AutoFit()
